$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1993865030674846
$ws.Range("C2").Value = 0.549079754601227
$ws.Range("J2").Value = 0.009202453987730062
$ws.Range("P2").Value = 0.1257668711656442
$ws.Range("S2").Value = 0.1165644171779141
$ws.Range("B3").Value = 0.01036269430051814
$ws.Range("C3").Value = 0.0310880829015544
$ws.Range("J3").Value = 0.0155440414507772
$ws.Range("P3").Value = 0.772020725388601
$ws.Range("S3").Value = 0.1709844559585492
$ws.Range("J4").Value = 0.03508771929824561
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.2807017543859649
$ws.Range("B6").Value = 0.05714285714285714
$ws.Range("D6").Value = 0.0163265306122449
$ws.Range("F6").Value = 0.0653061224489796
$ws.Range("J6").Value = 0.273469387755102
$ws.Range("O6").Value = 0.0326530612244898
$ws.Range("Q6").Value = 0.1428571428571428
$ws.Range("R6").Value = 0.08979591836734693
$ws.Range("S6").Value = 0.3224489795918367
$ws.Range("B7").Value = 0.1407035175879397
$ws.Range("D7").Value = 0.01507537688442211
$ws.Range("F7").Value = 0.07537688442211055
$ws.Range("J7").Value = 0.1457286432160804
$ws.Range("O7").Value = 0.03517587939698492
$ws.Range("Q7").Value = 0.1758793969849246
$ws.Range("R7").Value = 0.06532663316582915
$ws.Range("S7").Value = 0.3467336683417085
$ws.Range("B8").Value = 0.09567198177676538
$ws.Range("D8").Value = 0.01594533029612756
$ws.Range("E8").Value = 0.004555808656036446
$ws.Range("F8").Value = 0.05011389521640091
$ws.Range("J8").Value = 0.09339407744874716
$ws.Range("O8").Value = 0.03644646924829157
$ws.Range("Q8").Value = 0.2004555808656036
$ws.Range("R8").Value = 0.08200455580865604
$ws.Range("S8").Value = 0.4214123006833713
$ws.Range("B9").Value = 0.1065573770491803
$ws.Range("D9").Value = 0.02459016393442623
$ws.Range("F9").Value = 0.04098360655737705
$ws.Range("J9").Value = 0.06557377049180328
$ws.Range("O9").Value = 0.02868852459016394
$ws.Range("Q9").Value = 0.1844262295081967
$ws.Range("R9").Value = 0.0860655737704918
$ws.Range("S9").Value = 0.4631147540983607
$ws.Range("B10").Value = 0.1056081573197378
$ws.Range("D10").Value = 0.0269482884195193
$ws.Range("E10").Value = 0.001456664238892935
$ws.Range("F10").Value = 0.07137654770575383
$ws.Range("J10").Value = 0.101238164603059
$ws.Range("O10").Value = 0.0269482884195193
$ws.Range("Q10").Value = 0.2286962855061908
$ws.Range("R10").Value = 0.0764748725418791
$ws.Range("S10").Value = 0.3612527312454479
$ws.Range("F11").Value = 0.003205128205128205
$ws.Range("G11").Value = 0.1602564102564103
$ws.Range("J11").Value = 0.09615384615384616
$ws.Range("K11").Value = 0.2019230769230769
$ws.Range("L11").Value = 0.5256410256410257
$ws.Range("S11").Value = 0.01282051282051282
$ws.Range("G12").Value = 0.7108433734939759
$ws.Range("J12").Value = 0.2108433734939759
$ws.Range("K12").Value = 0.006024096385542169
$ws.Range("L12").Value = 0.03012048192771084
$ws.Range("S12").Value = 0.04216867469879518
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.2692307692307692
$ws.Range("S13").Value = 0.01923076923076923
$ws.Range("F15").Value = 0.01824817518248175
$ws.Range("H15").Value = 0.1277372262773723
$ws.Range("I15").Value = 0.06569343065693431
$ws.Range("J15").Value = 0.2992700729927008
$ws.Range("K15").Value = 0.0583941605839416
$ws.Range("M15").Value = 0.0145985401459854
$ws.Range("N15").Value = 0.0072992700729927
$ws.Range("O15").Value = 0.06204379562043796
$ws.Range("S15").Value = 0.3467153284671533
$ws.Range("F16").Value = 0.01310043668122271
$ws.Range("H16").Value = 0.148471615720524
$ws.Range("I16").Value = 0.1004366812227074
$ws.Range("J16").Value = 0.4104803493449782
$ws.Range("K16").Value = 0.1048034934497817
$ws.Range("M16").Value = 0.03056768558951965
$ws.Range("O16").Value = 0.06986899563318777
$ws.Range("S16").Value = 0.1222707423580786
$ws.Range("F17").Value = 0.02509652509652509
$ws.Range("H17").Value = 0.138996138996139
$ws.Range("I17").Value = 0.1042471042471042
$ws.Range("J17").Value = 0.4092664092664093
$ws.Range("K17").Value = 0.0945945945945946
$ws.Range("M17").Value = 0.01351351351351351
$ws.Range("O17").Value = 0.03667953667953668
$ws.Range("S17").Value = 0.1776061776061776
$ws.Range("F18").Value = 0.01507537688442211
$ws.Range("H18").Value = 0.1708542713567839
$ws.Range("I18").Value = 0.08542713567839195
$ws.Range("J18").Value = 0.4170854271356784
$ws.Range("K18").Value = 0.09547738693467336
$ws.Range("M18").Value = 0.02512562814070352
$ws.Range("O18").Value = 0.07035175879396985
$ws.Range("S18").Value = 0.1206030150753769
$ws.Range("F19").Value = 0.02410468319559229
$ws.Range("H19").Value = 0.1880165289256198
$ws.Range("I19").Value = 0.09297520661157024
$ws.Range("J19").Value = 0.3739669421487603
$ws.Range("K19").Value = 0.09297520661157024
$ws.Range("M19").Value = 0.02410468319559229
$ws.Range("O19").Value = 0.07024793388429752
$ws.Range("S19").Value = 0.1356749311294766
